$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value2 = 352.35715
$ws.Range("I12").Value2 = 286.16666
$ws.Range("J12").Value2 = 749.5
$ws.Range("K12").Value2 = 286.16666
$ws.Range("L12").Value2 = 749.5
$ws.Range("M12").Value2 = -116.16666
$ws.Range("N12").Value2 = -1089.5
$ws.Range("H70").Value2 = 826.9259
$ws.Range("I70").Value2 = 605.3333
$ws.Range("J70").Value2 = 1270.1111
$ws.Range("K70").Value2 = 1815.9999
$ws.Range("L70").Value2 = 3810.3333
$ws.Range("M70").Value2 = -1545.9999
$ws.Range("N70").Value2 = -4350.3333
$ws.Range("H73").Value2 = 826.9259
$ws.Range("I73").Value2 = 605.3333
$ws.Range("J73").Value2 = 1270.1111
$ws.Range("K73").Value2 = 1815.9999
$ws.Range("L73").Value2 = 3810.3333
$ws.Range("M73").Value2 = -879.9999
$ws.Range("N73").Value2 = -5682.3333
$ws.Range("H76").Value2 = 5128.4287
$ws.Range("I76").Value2 = 5316.5
$ws.Range("K76").Value2 = 5316.5
$ws.Range("M76").Value2 = -5001.5
$ws.Range("H79").Value2 = 5128.4287
$ws.Range("I79").Value2 = 5316.5
$ws.Range("K79").Value2 = 5316.5
$ws.Range("M79").Value2 = -4224.5
$ws.Range("H106").Value2 = 4582
$ws.Range("I106").Value2 = 4184.909
$ws.Range("K106").Value2 = 4184.909
$ws.Range("M106").Value2 = -3553.909
$ws.Range("H115").Value2 = 637.8570999999999
$ws.Range("I115").Value2 = 637.8570999999999
$ws.Range("K115").Value2 = 1913.5713
$ws.Range("M115").Value2 = -346.5712999999998
$ws.Range("H116").Value2 = 8455.1
$ws.Range("I116").Value2 = 7402.222
$ws.Range("K116").Value2 = 7402.222
$ws.Range("M116").Value2 = -3960.222
$ws.Range("H130").Value2 = 26664.834
$ws.Range("J130").Value2 = 26664.834
$ws.Range("L130").Value2 = 26664.834
$ws.Range("N130").Value2 = -36704.834
$ws.Range("H137").Value2 = 4135.1665
$ws.Range("J137").Value2 = 3329.5
$ws.Range("L137").Value2 = 9988.5
$ws.Range("N137").Value2 = -15088.5
$ws.Range("H138").Value2 = 2776.5
$ws.Range("J138").Value2 = 3638.2646
$ws.Range("L138").Value2 = 10914.7938
$ws.Range("N138").Value2 = -21194.7938

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 14380.885
$ws.Range("I2").Value2 = 19128.166
$ws.Range("J2").Value2 = 3699.5
$ws.Range("K2").Value2 = 19128.166
$ws.Range("L2").Value2 = 3699.5
$ws.Range("M2").Value2 = -19015.166
$ws.Range("N2").Value2 = -3925.5
$ws.Range("H116").Value2 = 14380.885
$ws.Range("I116").Value2 = 19128.166
$ws.Range("J116").Value2 = 3699.5
$ws.Range("K116").Value2 = 19128.166
$ws.Range("L116").Value2 = 3699.5
$ws.Range("M116").Value2 = -16834.166
$ws.Range("N116").Value2 = -8287.5
$ws.Range("H122").Value2 = 3673.4893
$ws.Range("I122").Value2 = 3391.1282
$ws.Range("K122").Value2 = 10173.3846
$ws.Range("M122").Value2 = -7723.384600000001
$ws.Range("H132").Value2 = 2809.639
$ws.Range("I132").Value2 = 2782.4517
$ws.Range("J132").Value2 = 2978.2
$ws.Range("K132").Value2 = 8347.355100000001
$ws.Range("L132").Value2 = 8934.599999999999
$ws.Range("M132").Value2 = -5817.355100000001
$ws.Range("N132").Value2 = -13994.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 14380.885
$ws.Range("I3").Value2 = 19128.166
$ws.Range("J3").Value2 = 3699.5
$ws.Range("K3").Value2 = 19128.166
$ws.Range("L3").Value2 = 3699.5
$ws.Range("M3").Value2 = -19014.166
$ws.Range("N3").Value2 = -3927.5
$ws.Range("H94").Value2 = 2055.8076
$ws.Range("I94").Value2 = 922.6
$ws.Range("K94").Value2 = 922.6
$ws.Range("M94").Value2 = -471.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1986.8
$ws.Range("I16").Value2 = 2063.6
$ws.Range("J16").Value2 = 1910
$ws.Range("K16").Value2 = 2063.6
$ws.Range("L16").Value2 = 1910
$ws.Range("M16").Value2 = -1776.6
$ws.Range("N16").Value2 = -2484
$ws.Range("H22").Value2 = 699.8333
$ws.Range("I22").Value2 = 733.3333
$ws.Range("K22").Value2 = 733.3333
$ws.Range("M22").Value2 = -383.3333
$ws.Range("H31").Value2 = 1615.5714
$ws.Range("I31").Value2 = 1460.4359
$ws.Range("J31").Value2 = 3632.3333
$ws.Range("K31").Value2 = 1460.4359
$ws.Range("L31").Value2 = 3632.3333
$ws.Range("M31").Value2 = -1165.4359
$ws.Range("N31").Value2 = -4222.3333
$ws.Range("H34").Value2 = 1615.5714
$ws.Range("I34").Value2 = 1460.4359
$ws.Range("J34").Value2 = 3632.3333
$ws.Range("K34").Value2 = 1460.4359
$ws.Range("L34").Value2 = 3632.3333
$ws.Range("M34").Value2 = -1258.4359
$ws.Range("N34").Value2 = -4036.3333
$ws.Range("H35").Value2 = 341.66666
$ws.Range("I35").Value2 = 341.66666
$ws.Range("K35").Value2 = 341.66666
$ws.Range("M35").Value2 = -47.66665999999998
$ws.Range("H86").Value2 = 32444.8
$ws.Range("I86").Value2 = 38199.89
$ws.Range("K86").Value2 = 38199.89
$ws.Range("M86").Value2 = -37076.89
$ws.Range("H89").Value2 = 32444.8
$ws.Range("I89").Value2 = 38199.89
$ws.Range("K89").Value2 = 190999.45
$ws.Range("M89").Value2 = -185383.45
$ws.Range("H113").Value2 = 1986.8
$ws.Range("I113").Value2 = 2063.6
$ws.Range("J113").Value2 = 1910
$ws.Range("K113").Value2 = 2063.6
$ws.Range("L113").Value2 = 1910
$ws.Range("M113").Value2 = 106.4000000000001
$ws.Range("N113").Value2 = -6250
$ws.Range("H132").Value2 = 2010.0646
$ws.Range("I132").Value2 = 1501.0385
$ws.Range("K132").Value2 = 4503.1155
$ws.Range("M132").Value2 = -1973.1155
$ws.Range("H134").Value2 = 1486.8
$ws.Range("I134").Value2 = 1233.5
$ws.Range("J134").Value2 = 2500
$ws.Range("K134").Value2 = 3700.5
$ws.Range("L134").Value2 = 7500
$ws.Range("M134").Value2 = -1165.5
$ws.Range("N134").Value2 = -12570

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value2 = 14231.25
$ws.Range("J12").Value2 = 18325
$ws.Range("L12").Value2 = 54975
$ws.Range("N12").Value2 = -55321
$ws.Range("H13").Value2 = 245
$ws.Range("I13").Value2 = 245
$ws.Range("J13").Value2 = 0
$ws.Range("K13").Value2 = 735
$ws.Range("L13").Value2 = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value2 = -567
$ws.Range("H16").Value2 = 2000
$ws.Range("I16").Value2 = 1000
$ws.Range("K16").Value2 = 3000
$ws.Range("M16").Value2 = -2827
$ws.Range("H20").Value2 = 1000
$ws.Range("I20").Value2 = 0
$ws.Range("J20").Value2 = 1000
$ws.Range("K20").Value2 = 0
$ws.Range("L20").ClearContents()
$ws.Range("M20").Value2 = 3000
$ws.Range("N20").Value2 = -3454
$ws.Range("H22").Value2 = 111111460
$ws.Range("I22").Value2 = 515.5
$ws.Range("J22").Value2 = 333333340
$ws.Range("K22").Value2 = 1546.5
$ws.Range("L22").Value2 = 1000000020
$ws.Range("M22").Value2 = -1377.5
$ws.Range("N22").Value2 = -1000000358
$ws.Range("H27").Value2 = 111111460
$ws.Range("I27").Value2 = 515.5
$ws.Range("J27").Value2 = 333333340
$ws.Range("K27").Value2 = 1546.5
$ws.Range("L27").Value2 = 1000000020
$ws.Range("M27").Value2 = -1444.5
$ws.Range("N27").Value2 = -1000000224
$ws.Range("H86").Value2 = 1235.5834
$ws.Range("I86").Value2 = 1576.3334
$ws.Range("J86").Value2 = 894.8333
$ws.Range("K86").Value2 = 4729.0002
$ws.Range("L86").Value2 = 2684.4999
$ws.Range("M86").Value2 = -3543.0002
$ws.Range("N86").Value2 = -5056.4999
$ws.Range("H89").Value2 = 1235.5834
$ws.Range("I89").Value2 = 1576.3334
$ws.Range("J89").Value2 = 894.8333
$ws.Range("K89").Value2 = 14187.0006
$ws.Range("L89").Value2 = 8053.4997
$ws.Range("M89").Value2 = -8259.000599999999
$ws.Range("N89").Value2 = -19909.4997
$ws.Range("H117").Value2 = 1174.6
$ws.Range("I117").Value2 = 511.4
$ws.Range("J117").Value2 = 1837.8
$ws.Range("K117").Value2 = 1534.2
$ws.Range("L117").Value2 = 5513.4
$ws.Range("M117").Value2 = 1907.8
$ws.Range("N117").Value2 = -12397.4
$ws.Range("H131").Value2 = 3317.8125
$ws.Range("I131").Value2 = 1719.25
$ws.Range("J131").Value2 = 4916.375
$ws.Range("K131").Value2 = 5157.75
$ws.Range("L131").Value2 = 14749.125
$ws.Range("M131").Value2 = -117.75
$ws.Range("N131").Value2 = -24829.125
$ws.Range("H132").Value2 = 2885.111
$ws.Range("I132").Value2 = 1989
$ws.Range("K132").Value2 = 17901
$ws.Range("M132").Value2 = -15371

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value2 = 0
$ws.Range("J47").Value2 = 0
$ws.Range("L47").ClearContents()
$ws.Range("N47").Value2 = 0
$ws.Range("H113").Value2 = 7353.778
$ws.Range("J113").Value2 = 9713
$ws.Range("L113").Value2 = 9713
$ws.Range("N113").Value2 = -14053
$ws.Range("H122").Value2 = 3781
$ws.Range("I122").Value2 = 2994.5
$ws.Range("K122").Value2 = 8983.5
$ws.Range("M122").Value2 = -6533.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value2 = 27799
$ws.Range("I42").Value2 = 19000
$ws.Range("J42").Value2 = 29998.75
$ws.Range("K42").Value2 = 19000
$ws.Range("L42").Value2 = 29998.75
$ws.Range("M42").Value2 = -18437
$ws.Range("N42").Value2 = -31124.75
$ws.Range("H49").Value2 = 27799
$ws.Range("I49").Value2 = 19000
$ws.Range("J49").Value2 = 29998.75
$ws.Range("K49").Value2 = 19000
$ws.Range("L49").Value2 = 29998.75
$ws.Range("M49").Value2 = -18853
$ws.Range("N49").Value2 = -30292.75
$ws.Range("H61").Value2 = 2758
$ws.Range("I61").Value2 = 2575.75
$ws.Range("J61").Value2 = 3001
$ws.Range("K61").Value2 = 2575.75
$ws.Range("L61").Value2 = 3001
$ws.Range("M61").Value2 = -2373.75
$ws.Range("N61").Value2 = -3405
$ws.Range("H93").Value2 = 2369.1667
$ws.Range("I93").Value2 = 2216.6
$ws.Range("J93").Value2 = 2559.875
$ws.Range("K93").Value2 = 2216.6
$ws.Range("L93").Value2 = 2559.875
$ws.Range("M93").Value2 = -968.5999999999999
$ws.Range("N93").Value2 = -5055.875
$ws.Range("H113").Value2 = 2758
$ws.Range("I113").Value2 = 2575.75
$ws.Range("J113").Value2 = 3001
$ws.Range("K113").Value2 = 2575.75
$ws.Range("L113").Value2 = 3001
$ws.Range("M113").Value2 = -405.75
$ws.Range("N113").Value2 = -7341

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value2 = 66308
$ws.Range("I96").Value2 = 128489.5
$ws.Range("J96").Value2 = 4126.5
$ws.Range("K96").Value2 = 128489.5
$ws.Range("L96").Value2 = 4126.5
$ws.Range("M96").Value2 = -127116.5
$ws.Range("N96").Value2 = -6872.5
$ws.Range("H107").Value2 = 362.2353
$ws.Range("I107").Value2 = 373.6875
$ws.Range("J107").Value2 = 179
$ws.Range("K107").Value2 = 1121.0625
$ws.Range("L107").Value2 = 537
$ws.Range("M107").Value2 = 798.9375
$ws.Range("N107").Value2 = -4377
